{"js": "// The document's placeholder demo-link line is updated to include the\n// actual Demo Video link, e.g.:\n//   (Add your demo link here)\n// becomes\n//   Demo Link: https://drive.google.com/file/d/1zmtP-gMoP59yaZZNCyndsJOaZmKvF0oo/view?usp=drive_link\nconst body = context.document.body;\nconst results = body.search(\"(Add your demo link here)\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Demo Link: https://drive.google.com/file/d/1zmtP-gMoP59yaZZNCyndsJOaZmKvF0oo/view?usp=drive_link\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# The placeholder demo-link line is updated to include the actual Demo\n# Video link, e.g.:\n#   (Add your demo link here)\n# becomes\n#   Demo Link: https://drive.google.com/file/d/1zmtP-gMoP59yaZZNCyndsJOaZmKvF0oo/view?usp=drive_link\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"(Add your demo link here)\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Demo Link: https://drive.google.com/file/d/1zmtP-gMoP59yaZZNCyndsJOaZmKvF0oo/view?usp=drive_link\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
